$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad / changed-date) for rows 2-34 from 45224 to 45233
$ws.Range("C2:C34").Value = 45233
